# Append a new row of scraped data to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on every existing data row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-17 12:53:41"

# Update the fetch timestamp on the already-existing rows (2-10).
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Append the new row (row 11).
$row = 11
$ws.Cells.Item($row, 1).Value = $newTimestamp
$ws.Cells.Item($row, 2).Value = "グループ毎の日報をリアルタイムでまとめたい"
$ws.Cells.Item($row, 3).Value = "システム開発"
$ws.Cells.Item($row, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item($row, 5).Value = "期限情報なし"

$url = "https://www.lancers.jp/work/detail/5456195"
$ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $url)
$ws.Cells.Item($row, 6).Style = "Hyperlink"

$ws.Cells.Item($row, 7).Value = 18
